$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 13 (old rows 13-23 shift down to 14-24)
$ws.Rows.Item(13).Insert()

# 2) The inserted row 13 picked up a copied A13 cell (bold style) from row above; remove it
$ws.Cells.Item(13, 1).Clear()

# 3) Copy column B/C formatting from row 10 (Objetivos row) onto row 13, then set its text
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4122)
$ws.Cells.Item(10, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial(-4122)
$ws.Cells.Item(13, 2).Value = '11079086 - Herlandí de Souza Andrade'
$ws.Cells.Item(13, 3).Value = '11079086 - Herlandí de Souza Andrade'

# 4) Update "Objetivos:" row (row 10) B/C text
$ws.Cells.Item(10, 2).Value = 'Oferecer ao aluno uma visão geral dos conceitos fundamentais e das fases do processo de desenvolvimento e Implementação de Sistemas de Informação no sentido de capacitá-lo analisa e projetar tais sistemas'
$ws.Cells.Item(10, 3).Value = 'Oferecer ao aluno uma visão geral dos conceitos fundamentais e das fases do processo de desenvolvimento e Implementação de Sistemas de Informação no sentido de capacitá-lo analisa e projetar tais sistemas'

# 5) Update "Programa resumido:" row (now row 14) B/C text
$ws.Cells.Item(14, 2).Value = 'Sistemas de Informação. Projeto de Sistemas de Informação. Tecnologia de Informação. Processo de Desenvolvimento de Sistema de Informação.'
$ws.Cells.Item(14, 3).Value = 'Sistemas de Informação. Projeto de Sistemas de Informação. Tecnologia de Informação. Processo de Desenvolvimento de Sistema de Informação.'

# 6) Update "Programa:" row (now row 16) B/C text
$ws.Cells.Item(16, 2).Value = '1. Sistemas de Informação1.1. Sistemas de Processamento de Informações;1.2. Sistemas de Informações Gerenciais;1.3. Sistema de Apoio à Decisão;1.4. Sistemas de Informação no Comércio Eletrônico;1.5. Sistemas de Informação em Cadeia de Suprimentos;1.6. Sistemas inteligentes nos negócios;1.7. Sistemas estratégicos. 2. Projeto de Sistemas de Informação.2.1. Especificação das Saídas;2.2. Especificação dos Arquivos;2.3. Especificação das Entradas;2.4. Especificação do Processamento.3. Tecnologia de Informação.3.1. Evolução da Computação;3.2. Recursos Computacionais.4. Processo de Desenvolvimento de Sistemas de Informação.4.1. Definição do Negócio;4.2. Identificação do Problema e/ou Oportunidades;4.3. Seleção do Sistema de Informação;4.4. Implementação do Sistema de Informação;4.5. Avaliação da Eficácia do Sistema de Informação;'
$ws.Cells.Item(16, 3).Value = '1. Sistemas de Informação1.1. Sistemas de Processamento de Informações;1.2. Sistemas de Informações Gerenciais;1.3. Sistema de Apoio à Decisão;1.4. Sistemas de Informação no Comércio Eletrônico;1.5. Sistemas de Informação em Cadeia de Suprimentos;1.6. Sistemas inteligentes nos negócios;1.7. Sistemas estratégicos. 2. Projeto de Sistemas de Informação.2.1. Especificação das Saídas;2.2. Especificação dos Arquivos;2.3. Especificação das Entradas;2.4. Especificação do Processamento.3. Tecnologia de Informação.3.1. Evolução da Computação;3.2. Recursos Computacionais.4. Processo de Desenvolvimento de Sistemas de Informação.4.1. Definição do Negócio;4.2. Identificação do Problema e/ou Oportunidades;4.3. Seleção do Sistema de Informação;4.4. Implementação do Sistema de Informação;4.5. Avaliação da Eficácia do Sistema de Informação;'

# 7) Update "Método:" row (now row 19) B/C text
$ws.Cells.Item(19, 2).Value = 'Aulas expositivas teóricas, aulas práticas, aulas de exercícios.'
$ws.Cells.Item(19, 3).Value = 'Aulas expositivas teóricas, aulas práticas, aulas de exercícios.'

# 8) Update "Critério:" row (now row 20) B/C text
$ws.Cells.Item(20, 2).Value = 'Média Aritmética das atividades avaliativas realizadas.'
$ws.Cells.Item(20, 3).Value = 'Média Aritmética das atividades avaliativas realizadas.'

# 9) Update "Norma de recuperação:" row (now row 21) B/C text
$ws.Cells.Item(21, 2).Value = 'Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.'
$ws.Cells.Item(21, 3).Value = 'Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.'

# 10) Update "Bibliografia:" row (now row 22) B/C text
$ws.Cells.Item(22, 2).Value = 'HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004LAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gestão integrada de processos e da tecnologia da informação. São Paulo:Atlas, 2006.LAURINDO, F.J.B. Tecnologia da Informação: Eficácia nas Organizações. São Paulo, Editora Futura, 2002.STAIR, R.M., Princípios de Sistema de Informação: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.TURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.TURBAN, E., RAIANER JR, K., POTTER, R. E., Administração de Tecnologia da Informação: Teoria e Prática”, São Paulo, Editora Campus, 2003.'
$ws.Cells.Item(22, 3).Value = 'HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004LAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gestão integrada de processos e da tecnologia da informação. São Paulo:Atlas, 2006.LAURINDO, F.J.B. Tecnologia da Informação: Eficácia nas Organizações. São Paulo, Editora Futura, 2002.STAIR, R.M., Princípios de Sistema de Informação: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.TURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.TURBAN, E., RAIANER JR, K., POTTER, R. E., Administração de Tecnologia da Informação: Teoria e Prática”, São Paulo, Editora Campus, 2003.'

# 11) Fix the column definitions: column A ("min=1 max=2") used to share a width
#     entry with column B. Touch column B's width so the engine splits the
#     redundant combined range into a distinct column A entry.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
